$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily observations appended after row 189 (30-09-2021).
# Force column A to text first so the dd-mm-yyyy strings are stored as
# shared strings instead of being auto-converted to date serials, then
# drop the temporary number format so the cells keep the sheet's default
# (unstyled) look, matching the existing rows above them.
$ws.Range("A190:A191").NumberFormat = "@"

$ws.Range("A190").Value = "01-10-2021"
$ws.Range("B190").Value = 1.73
$ws.Range("C190").Value = 2.49
$ws.Range("D190").Value = 3.2
$ws.Range("E190").Value = 3.84
$ws.Range("F190").Value = -0.18

$ws.Range("A191").Value = "04-10-2021"
$ws.Range("B191").Value = 1.72
$ws.Range("C191").Value = 2.51
$ws.Range("D191").Value = 3.23
$ws.Range("E191").Value = 3.88
$ws.Range("F191").Value = -0.17

$ws.Range("A190:A191").Style = "Normal"
